$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Naive-forecaster bugfix: a previously-missing "current quarter" error
# value now gets written into column B of each data row; the values that
# used to start at column B shift one column to the right (B->C, C->D,
# ...), with whatever falls past the row's original last column dropped.
# (Numbers are written in plain-decimal form -- this engine's script
# parser doesn't accept scientific-notation literals.)
$newValues = @{
    2  = -0.0000003965936795080616225
    3  = -0.0000000003930720193778825
    4  = 0.0000001181302580199883058
    5  = 0.0000002965444589886345739
    6  = 0.0000003593882045849206008
    7  = -0.0000001035781544145298037
    8  = -0.000000000387512216759589
    9  = -0.0000001073887893610070024
    10 = -0.000000103547280583260499
    11 = 0.0000063033553409086451325
    12 = -0.0000002375649628613696064
    13 = 0.0000003720025918141355884
    14 = 0.0000003829984367986761203
    15 = -0.0000031604754923975080033
    16 = -0.0000000410109615434084418
    17 = -0.0000001831659499074156107
    18 = 0.0000002770877186031305865
    19 = 0.0000002297750048008140046
    20 = -0.00000015542410669588949
}

for ($row = 2; $row -le 20; $row++) {

    # Read the existing values starting at column B (2) until the first
    # empty cell in the row.
    $oldValues = @()
    $col = 2
    while ($true) {
        $cell = $ws.Cells.Item($row, $col)
        $v = $cell.Value()
        if ($v -eq $null) { break }
        $oldValues += , $v
        $col++
    }

    # Shift the captured values one column to the right (dropping the
    # last one, which falls outside the row's original extent), then
    # write the new value into column B.
    if ($oldValues.Count -gt 0) {
        for ($i = 0; $i -lt ($oldValues.Count - 1); $i++) {
            $ws.Cells.Item($row, 3 + $i).Value = $oldValues[$i]
        }
    }

    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}
